# Fixed tests post changes to Stakeholder
# - Rename the "Investor" header (A1) to "Stakeholder"
# - Update the active selection to A2 (was K4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Stakeholder"

[void]$ws.Range("A2").Select()
